# Updates the localized text strings for the "language" workbook and
# adjusts the remembered sheet view/selection (level 4 section 3 play).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Build a lookup of Key text (column A) -> row number so we can find rows
# by their key name instead of relying on fixed row numbers.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$keyToRow = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $k = $ws.Cells.Item($r, 1).Value2
    if ($k -ne $null) {
        $keyToRow[[string]$k] = $r
    }
}

function Set-ValueForKey($Key, $Text) {
    $row = $keyToRow[$Key]
    if ($row -eq $null) {
        throw "Key not found: $Key"
    }
    $ws.Cells.Item($row, 2).Value2 = $Text
}

$updates = @(
    @("level01_intro_1", "To help the mighty Chin-Chin, we need to learn how to use objects with special properties."),
    @("level01_intro_2", "Each object has properties that fit different tasks."),
    @("level01_property_1", "Objects have many different properties."),
    @("level01_property_2", "Each property tells us its quality and attribute."),
    @("level01_property_3", "For example, a rubber duck is: lightweight, not very strong, floats on water, and can block electricity."),
    @("level01_classify_1", "Classify means putting objects into groups with shared properties."),
    @("level01_classify_2", "This allows us to use the right object for every task!"),
    @("level01_shape_size_1", "Objects can be different shapes and sizes."),
    @("level01_begin_1", "Let's go ahead and shape clay!"),
    @("level01_begin_2", "Place the clay into two groups: tall and wide."),
    @("level01_play_begin_1", "Look! There's a star fragment."),
    @("level01_play_begin_2", "Let's help the mighty Chin-Chin cross the deadly pit!"),
    @("level02_intro_1", "Now we will look at the weight property of objects."),
    @("level02_intro_2", "An object's weight tells us how light or heavy it is."),
    @("level02_intro_3", "You can determine an object's weight with a balance scale."),
    @("level02_light_1", "Here, a 2-pound pillow is light compared to a 20-pound block."),
    @("level02_heavy_1", "A 40-pound iron block is heavycompared to a 20-pound block."),
    @("level02_begin_1", "Now let's sort objects into two classifications: light and heavy!"),
    @("level02_begin_2", "An object is light if the left scale goes up, and heavy if the left scale goes down."),
    @("level03_intro_1", "Now let's look at buoyancy."),
    @("level03_buoyancy_1", "An object's buoyancy determines whether it floats or sinks in water."),
    @("level03_buoyancy_2", "If an object moves up in water, it floats."),
    @("level03_buoyancy_3", "If an object moves down in water, it sinks."),
    @("level03_begin_1", "Sort the objects on whether they float or sink!"),
    @("level04_intro_1", "The final mission! We learn about the conductive property."),
    @("level04_conductive_1", "The iron block is highly conductive with electricity. It lets electricity flow through it very easily."),
    @("level04_non_conductive_1", "The granite rock is not conductive."),
    @("level04_begin_1", "Now, let's go ahead and sort these objects on whether they are electrically conductive or not!"),
    @("level04_begin_2", "Sort these objects by electrically conductive or not!"),
    @("end_title", "CONGRATULATIONS!")
)

foreach ($u in $updates) {
    Set-ValueForKey $u[0] $u[1]
}

# Restore the saved view/selection state to match the author's final edit
# position (scrolled further down the sheet, selection moved to B90).
$ws.Activate()
$ws.Range("B90").Select()
